# Update for March 18 — append the new day's row (row 24) to the
# OntarioCoronavirus.csv sheet, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OntarioCoronavirus.csv")

$newUrl = "https://www.cbc.ca/news/canada/toronto/coronavirus-covid-19-ontario-wednesday-1.5501250"

# --- 1. New row of raw data -------------------------------------------------
$ws.Range("A24").Value = 43908
$ws.Range("B24").Value = 212
$ws.Range("C24").Value = 13897

# Formulas mirror the shared-formula pattern already used down column D/E/F.
$ws.Range("D24").Formula = '=A24-$A$2'
$ws.Range("E24").Formula = '=D24-D23'
$ws.Range("F24").Formula = '=(B24/B23)^(1/E24)-1'

$ws.Range("G24").Value = $newUrl

# --- 2. Match formatting of the row above (date style, percent style, ...) -
$ws.Range("A23:G23").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)  # xlPasteFormats

# --- 3. Hyperlink for the new source cell -----------------------------------
$ws.Hyperlinks.Add($ws.Range("G24"), $newUrl)
$ws.Range("G24").Style = "Hyperlink"

# --- 4. Nudge the chart down by one row (it sits right below the table) ----
$co = $ws.ChartObjects(1)
$co.Top = $co.Top + 15

$excel.CutCopyMode = $false
